$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to text so numeric-looking strings
# (e.g. "23.70", "66.512.71") are not reinterpreted as numbers/dates.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '66.512.71'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '3.468.31'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '586.12'
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('D6').Value = '177.83'
$ws.Range('E6').Value = '  +0.74%  '
$ws.Range('D7').Value = '0.627'
$ws.Range('E7').Value = '  +4.96%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '3.464.06'
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('D10').Value = '0.133'
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('D13').Value = '4.069.56'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('E14').Value = '  +1.48%  '
$ws.Range('D15').Value = '30.25'
$ws.Range('E15').Value = '  -0.30%  '
$ws.Range('D16').Value = '66.344.44'
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '3.466.97'
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('D19').Value = '5.98'
$ws.Range('E19').Value = '  -1.15%  '
$ws.Range('E20').Value = '  -1.12%  '
$ws.Range('D21').Value = '372.71'
$ws.Range('E21').Value = '  -2.47%  '
$ws.Range('E22').Value = '  -2.01%  '
$ws.Range('D23').Value = '73.38'
$ws.Range('E23').Value = '  +1.48%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').Value = '0.538'
$ws.Range('E25').Value = '  -1.85%  '
$ws.Range('E26').Value = '  +4.68%  '
$ws.Range('D27').Value = '10.06'
$ws.Range('E27').Value = '  +2.25%  '
$ws.Range('E28').Value = '  +3.05%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').Value = '5.98'
$ws.Range('E30').Value = '  +1.67%  '
$ws.Range('E31').Value = '  -0.59%  '
$ws.Range('D32').Value = '23.70'
$ws.Range('E32').Value = '  -3.53%  '
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('D34').Value = '7.05'
$ws.Range('E34').Value = '  -2.91%  '
$ws.Range('E35').Value = '  -4.99%  '
$ws.Range('E36').Value = '  -0.60%  '
$ws.Range('D37').Value = '161.14'
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').Value = '0.885'
$ws.Range('E38').Value = '  -0.79%  '
$ws.Range('D39').Value = '28.19'
$ws.Range('E39').Value = '  -4.98%  '
$ws.Range('D40').Value = '1.82'
$ws.Range('E40').Value = '  +1.46%  '
$ws.Range('D41').Value = '2.803.20'
$ws.Range('E41').Value = '  +2.73%  '
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '6.50'
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '2.56'
$ws.Range('E44').Value = '  +0.67%  '
$ws.Range('D45').Value = '0.0696'
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('D46').Value = '25.19'
$ws.Range('E46').Value = '  -0.17%  '
$ws.Range('D47').Value = '340.47'
$ws.Range('E47').Value = '  +4.56%  '
$ws.Range('D48').Value = '40.04'
$ws.Range('E48').Value = '  -1.44%  '
$ws.Range('D49').Value = '0.0293'
$ws.Range('E49').Value = '  +0.40%  '
$ws.Range('E50').Value = '  +2.44%  '
$ws.Range('B51').Value = 'Arweave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D51').Value = '31.77'
$ws.Range('E51').Value = '  +2.30%  '

# Restore default (unstyled) formatting on the touched range so it
# matches the original unstyled data cells.
$ws.Range('D2:E51').ClearFormats()
